# Updated cryptos list on Fri Aug 25 15:37:36 UTC 2023 with GitHub Actions
#
# Refreshes the coin price / 1h-volume figures pulled from coinranking.com,
# and reflects a ranking swap between WrappedEther and Polkadot (rows 12/13).
#
# Numeric-looking price strings (e.g. "216.61") are written with a leading
# apostrophe so Excel stores them as text -- matching how the Price column
# is already populated (plain display strings like "26.017.07", "4.461",
# not real numbers) instead of letting Excel auto-convert them to floats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '26.017.07' }
    @{ Cell = 'E2'; Value = '  -0.74%  ' }
    @{ Cell = 'D3'; Value = '1.651.93' }
    @{ Cell = 'E3'; Value = '  -0.11%  ' }
    @{ Cell = 'E4'; Value = '  -0.31%  ' }
    @{ Cell = 'D5'; Value = '''216.61' }
    @{ Cell = 'E5'; Value = '  -0.72%  ' }
    @{ Cell = 'D6'; Value = '''0.5210' }
    @{ Cell = 'E6'; Value = '  +0.28%  ' }
    @{ Cell = 'E7'; Value = '  -0.33%  ' }
    @{ Cell = 'D8'; Value = '''0.2614' }
    @{ Cell = 'E8'; Value = '  -1.25%  ' }
    @{ Cell = 'D9'; Value = '''0.06265' }
    @{ Cell = 'E9'; Value = '  -0.60%  ' }
    @{ Cell = 'D10'; Value = '''20.59' }
    @{ Cell = 'E10'; Value = '  -3.09%  ' }
    @{ Cell = 'D11'; Value = '''0.07737' }
    @{ Cell = 'E11'; Value = '  +0.03%  ' }
    @{ Cell = 'B12'; Value = 'WrappedEther' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' }
    @{ Cell = 'D12'; Value = '1.676.39' }
    @{ Cell = 'E12'; Value = '  +1.43%  ' }
    @{ Cell = 'B13'; Value = 'Polkadot' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' }
    @{ Cell = 'D13'; Value = '''4.461' }
    @{ Cell = 'E13'; Value = '  +0.93%  ' }
    @{ Cell = 'D14'; Value = '1.879.62' }
    @{ Cell = 'E14'; Value = '  -0.01%  ' }
    @{ Cell = 'D15'; Value = '''0.5419' }
    @{ Cell = 'E15'; Value = '  -0.51%  ' }
    @{ Cell = 'D16'; Value = '0.0₅8096' }
    @{ Cell = 'E16'; Value = '  -1.20%  ' }
    @{ Cell = 'D17'; Value = '''64.86' }
    @{ Cell = 'E17'; Value = '  +0.33%  ' }
    @{ Cell = 'D18'; Value = '26.035.25' }
    @{ Cell = 'E18'; Value = '  -0.68%  ' }
    @{ Cell = 'E19'; Value = '  -0.35%  ' }
    @{ Cell = 'D20'; Value = '''4.565' }
    @{ Cell = 'E20'; Value = '  -2.45%  ' }
    @{ Cell = 'D21'; Value = '''191.30' }
    @{ Cell = 'E21'; Value = '  +0.26%  ' }
    @{ Cell = 'D22'; Value = '''10.00' }
    @{ Cell = 'E22'; Value = '  -1.64%  ' }
    @{ Cell = 'D23'; Value = '''5.975' }
    @{ Cell = 'E23'; Value = '  -3.33%  ' }
    @{ Cell = 'E24'; Value = '  -0.38%  ' }
    @{ Cell = 'D25'; Value = '''138.17' }
    @{ Cell = 'E25'; Value = '  -0.51%  ' }
    @{ Cell = 'E26'; Value = '  -0.71%  ' }
    @{ Cell = 'D27'; Value = '''7.241' }
    @{ Cell = 'E27'; Value = '  -0.49%  ' }
    @{ Cell = 'E28'; Value = '  +0.25%  ' }
    @{ Cell = 'D29'; Value = '''1.395' }
    @{ Cell = 'E29'; Value = '  -1.50%  ' }
    @{ Cell = 'D30'; Value = '''0.05959' }
    @{ Cell = 'E30'; Value = '  -1.74%  ' }
    @{ Cell = 'D31'; Value = '''1.271' }
    @{ Cell = 'E31'; Value = '  -0.88%  ' }
    @{ Cell = 'D32'; Value = '''3.500' }
    @{ Cell = 'E32'; Value = '  -1.26%  ' }
    @{ Cell = 'D33'; Value = '''3.239' }
    @{ Cell = 'E33'; Value = '  -3.46%  ' }
    @{ Cell = 'D34'; Value = '''1.559' }
    @{ Cell = 'E34'; Value = '  -5.55%  ' }
    @{ Cell = 'D35'; Value = '''0.9482' }
    @{ Cell = 'E35'; Value = '  -3.63%  ' }
    @{ Cell = 'E36'; Value = '  -0.13%  ' }
    @{ Cell = 'D37'; Value = '''2.754' }
    @{ Cell = 'E37'; Value = '  -0.57%  ' }
    @{ Cell = 'D38'; Value = '''0.5679' }
    @{ Cell = 'E38'; Value = '  -4.20%  ' }
    @{ Cell = 'D39'; Value = '''0.01596' }
    @{ Cell = 'E39'; Value = '  +0.12%  ' }
    @{ Cell = 'D40'; Value = '''5.884' }
    @{ Cell = 'E40'; Value = '  -1.26%  ' }
    @{ Cell = 'D41'; Value = '''0.8452' }
    @{ Cell = 'E41'; Value = '  -1.98%  ' }
    @{ Cell = 'D42'; Value = '''1.001' }
    @{ Cell = 'E42'; Value = '  -0.29%  ' }
    @{ Cell = 'D43'; Value = '''100.74' }
    @{ Cell = 'E43'; Value = '  +0.96%  ' }
    @{ Cell = 'D44'; Value = '1.004.20' }
    @{ Cell = 'E44'; Value = '  -5.07%  ' }
    @{ Cell = 'D45'; Value = '1.795.12' }
    @{ Cell = 'E45'; Value = '  +0.03%  ' }
    @{ Cell = 'E46'; Value = '  -1.55%  ' }
    @{ Cell = 'D47'; Value = '''56.61' }
    @{ Cell = 'E47'; Value = '  -1.10%  ' }
    @{ Cell = 'D48'; Value = '''1.002' }
    @{ Cell = 'E48'; Value = '  -0.07%  ' }
    @{ Cell = 'D49'; Value = '''7.964' }
    @{ Cell = 'E49'; Value = '  -1.14%  ' }
    @{ Cell = 'D50'; Value = '''0.4302' }
    @{ Cell = 'E50'; Value = '  +1.62%  ' }
    @{ Cell = 'D51'; Value = '''1.473' }
    @{ Cell = 'E51'; Value = '  +0.24%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
